# ---------------------------------------------------------------------------
# RS-SAS-RDV.schema.docx edits
# ---------------------------------------------------------------------------
$d = $word.ActiveDocument

# 1. Heading: createAppointment -> appointment
$d.Content.Find.Execute("createAppointment", $true, $false, $false, $false, $false, $true, 1, $false, "appointment", 2) | Out-Null

# ---------------------------------------------------------------------------
# Table 1 ("appointment")
# ---------------------------------------------------------------------------
$t1 = $d.Tables(1)

# 2. Insert a new "method" row right after the "appointmentId" row (row 2),
#    i.e. before the current row 3 ("created").
$methodRow = $t1.Rows.Add($t1.Rows(3))
$methodRow.Cells(1).Range.Text = "method"
$methodRow.Cells(2).Range.Text = "Méthode"
$methodRow.Cells(3).Range.Text = "string|||(ENUM: CreateAppointment, UpdateAppointment)"
$methodRow.Cells(3).Range.Find.Execute("|||", $false, $false, $false, $false, $false, $true, 1, $false, "^l", 2) | Out-Null
$methodRow.Cells(4).Range.Text = "1..1"
$methodRow.Cells(5).Range.Text = "Indique un message de création ou de modification du rendez-vous"
$methodRow.Cells(6).Range.Text = "createAppointment"

# After the insertion the rows shift down by one:
#   3=method(new) 4=created 5=start 6=end 7=status 8=orientationCategory
#   9=practitioner 10=organization

# 3. "end" row cardinality: 1..1 -> 0..1
$t1.Cell(6,4).Range.Text = "0..1"

# 4. "status" row Format: add ENUM note
$t1.Cell(7,3).Range.Text = "string|||(ENUM: pending, booked, fulfilled, noshow, cancelled)"
$t1.Cell(7,3).Range.Find.Execute("|||", $false, $false, $false, $false, $false, $true, 1, $false, "^l", 2) | Out-Null

# 5. "orientationCategory" row Format: add ENUM note
$t1.Cell(8,3).Range.Text = "string|||(ENUM: CPTS, MSP, CDS, SOS, PS, PDM)"
$t1.Cell(8,3).Range.Find.Execute("|||", $false, $false, $false, $false, $false, $true, 1, $false, "^l", 2) | Out-Null

# 6. "orientationCategory" row Example: 604 -> SOS
$t1.Cell(8,6).Range.Text = "SOS"

# ---------------------------------------------------------------------------
# Table 2 ("practitioner")
# ---------------------------------------------------------------------------
$t2 = $d.Tables(2)

# 7. "rppsId" row Format: add REGEX note
$t2.Cell(2,3).Range.Text = "string|||(REGEX: ^81[0-9]{10}`$)"
$t2.Cell(2,3).Range.Find.Execute("|||", $false, $false, $false, $false, $false, $true, 1, $false, "^l", 2) | Out-Null

# 8. "firstName" row Cardinalite: 1..n -> 1..1
$t2.Cell(4,4).Range.Text = "1..1"

# 9. "speciality" -> "specialityCode"
$t2.Cell(5,1).Range.Text = "specialityCode"

# 10. "specialityCode" row Cardinalite: 0..n -> 0..1
$t2.Cell(5,4).Range.Text = "0..1"

# 11. Append three new rows after "specialityCode" (currently the last row).
$specialityUrlRow = $t2.Rows.Add()
$specialityUrlRow.Cells(1).Range.Text = "specialityUrl"
$specialityUrlRow.Cells(2).Range.Text = "Terminologie spécialité"
$specialityUrlRow.Cells(3).Range.Text = "string"
$specialityUrlRow.Cells(4).Range.Text = "0..1"
$specialityUrlRow.Cells(5).Range.Text = "Url de la terminologie utilisée pour la spécialité"
$specialityUrlRow.Cells(6).Range.Text = "https://mos.esante.gouv.fr/NOS/TRE_R38-SpecialiteOrdinale/FHIR/TRE-R38-SpecialiteOrdinale"

$professionCodeRow = $t2.Rows.Add()
$professionCodeRow.Cells(1).Range.Text = "professionCode"
$professionCodeRow.Cells(2).Range.Text = "Profession"
$professionCodeRow.Cells(3).Range.Text = "string"
$professionCodeRow.Cells(4).Range.Text = "0..1"
$professionCodeRow.Cells(5).Range.Text = "Code de la profession du professionnel de santé"
$professionCodeRow.Cells(6).Range.Text = "10"

$professionUrlRow = $t2.Rows.Add()
$professionUrlRow.Cells(1).Range.Text = "professionUrl"
$professionUrlRow.Cells(2).Range.Text = "Terminologie profession"
$professionUrlRow.Cells(3).Range.Text = "string"
$professionUrlRow.Cells(4).Range.Text = "0..1"
$professionUrlRow.Cells(5).Range.Text = "Url de la terminologie utilisée pour la profession"
$professionUrlRow.Cells(6).Range.Text = "https://mos.esante.gouv.fr/NOS/TRE_G15-ProfessionSante/FHIR/TRE-G15-ProfessionSante"

# ---------------------------------------------------------------------------
# Table 3 ("organization")
# ---------------------------------------------------------------------------
$t3 = $d.Tables(3)

# 12. "organizationId" row Exemple: 050005917 -> 334173748400020
$t3.Cell(2,6).Range.Text = "334173748400020"

# 13. "name" row Exemple: CDS DENTAIRE -> SOS Médecins de Rennes
$t3.Cell(3,6).Range.Text = "SOS Médecins de Rennes"

Write-Output "edits applied"
